$d = $word.ActiveDocument

# Change 1: "... an installer, simply run it by ..." -> "... an installer, run it by ..."
$d.Content.Find.Execute(
    "Once you have chosen and downloaded an installer, simply run it by double-clicking on the downloaded file. A dialog should appear that looks something like this:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Once you have chosen and downloaded an installer, run it by double-clicking on the downloaded file. A dialog should appear that looks something like this:",
    2
)

# Change 2: "Verify it" -> "Please verify it", ending up split across two runs
# ("Please v" / "erify it") exactly as the target XML has it.
$r = $d.Content
$r.Find.Execute("Verify it", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
$end = $r.End

# Insert "Please v" right before "Verify it" (same formatting, so it just
# extends the existing run for now).
$r.Collapse(1)
$r.InsertBefore("Please v")

# Remove the leading "V" of "Verify", turning it into "erify it".
$vChar = $d.Range($start + 8, $start + 9)
$vChar.Delete()

# Toggling a character-formatting property off/on forces the engine to
# materialize "erify it" as its own run even though the formatting ends up
# identical to its neighbour, matching the two-run split in the target doc.
# "Verify it" (9 chars) became "Please verify it" (17 chars): net +8 (the
# "Please v" insert) -1 (the deleted "V") = +7 versus the original $end.
$tail = $d.Range($start + 8, $end + 7)
$tail.Bold = 1
$tail.Bold = 0
